$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.324.33"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "3.792.03"
$ws.Range("E3").Value = "  +1.06%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'594.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'168.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").Value = "3.792.37"
$ws.Range("E7").Value = "  +1.14%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").Value = "'36.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.39%  "

$ws.Range("D15").Value = "4.426.65"
$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").Value = "3.793.52"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").Value = "68.353.41"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("E20").Value = "  -1.37%  "

$ws.Range("D21").Value = "'10.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "'465.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "'0.700"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("E24").Value = "  +9.04%  "

$ws.Range("D25").Value = "'83.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").Value = "'11.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("D28").Value = "'10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -1.02%  "

$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").Value = "'30.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "

$ws.Range("D33").Value = "'2.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "3.744.33"
$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").Value = "'3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("E40").Value = "  +0.89%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "'44.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.33%  "

$ws.Range("E45").Value = "  -2.39%  "

$ws.Range("D46").Value = "'47.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "

$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").Value = "'8.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("D49").Value = "'394.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").Value = "'146.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "

$ws.Range("D51").Value = "2.801.95"
$ws.Range("E51").Value = "  +4.29%  "
